# Commit: "error solve ifrs list"
# Fix the 2014-2021 IFRS financial figures on the company_list sheet for 신풍제약:
#  - rows 2-6 (FY2014-FY2018 actuals) get the corrected per-share-scaled figures
#  - rows 7-9 (FY2019E-FY2021E estimates) are cleared - they held bad/duplicated data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2014/12 (IFRS연결)
$ws.Range("D2").Value = 2203
$ws.Range("E2").Value = 193
$ws.Range("F2").Value = 193
$ws.Range("G2").Value = 135
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4093
$ws.Range("L2").Value = 2178
$ws.Range("M2").Value = 1915
$ws.Range("N2").Value = 1912
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 229
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -132
$ws.Range("S2").Value = 93
$ws.Range("T2").Value = 74
$ws.Range("U2").Value = -74
$ws.Range("V2").Value = 1570
$ws.Range("W2").Value = 8.74
$ws.Range("X2").Value = 1.52
$ws.Range("Y2").Value = 1.75
$ws.Range("Z2").Value = 0.83
$ws.Range("AA2").Value = 113.72
$ws.Range("AB2").Value = 805.35
$ws.Range("AC2").Value = 71
$ws.Range("AD2").Value = 60.73
$ws.Range("AE2").Value = 4530
$ws.Range("AF2").Value = 0.95
$ws.Range("AG2").Value = 48
$ws.Range("AH2").Value = 1.12
$ws.Range("AI2").Value = 61.47
$ws.Range("AJ2").Value = 45218000

# Row 3 - 2015/12 (IFRS연결)
$ws.Range("D3").Value = 1960
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 42
$ws.Range("G3").Value = 22
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4107
$ws.Range("L3").Value = 2159
$ws.Range("M3").Value = 1949
$ws.Range("N3").Value = 1946
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 237
$ws.Range("Q3").Value = 99
$ws.Range("R3").Value = -1
$ws.Range("S3").Value = -32
$ws.Range("T3").Value = 19
$ws.Range("U3").Value = 80
$ws.Range("V3").Value = 1561
$ws.Range("W3").Value = 2.16
$ws.Range("X3").Value = 0.75
$ws.Range("Y3").Value = 0.77
$ws.Range("Z3").Value = 0.36
$ws.Range("AA3").Value = 110.78
$ws.Range("AB3").Value = 775.0700000000001
$ws.Range("AC3").Value = 31
$ws.Range("AD3").Value = 140.05
$ws.Range("AE3").Value = 4610
$ws.Range("AF3").Value = 0.95
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 45218000

# Row 4 - 2016/12 (IFRS연결)
$ws.Range("D4").Value = 1941
$ws.Range("E4").Value = 96
$ws.Range("F4").Value = 96
$ws.Range("G4").Value = -37
$ws.Range("H4").Value = -186
$ws.Range("I4").Value = -186
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4165
$ws.Range("L4").Value = 1996
$ws.Range("M4").Value = 2168
$ws.Range("N4").Value = 2165
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 276
$ws.Range("Q4").Value = 190
$ws.Range("R4").Value = -212
$ws.Range("S4").Value = 55
$ws.Range("T4").Value = 14
$ws.Range("U4").Value = 175
$ws.Range("V4").Value = 1214
$ws.Range("W4").Value = 4.93
$ws.Range("X4").Value = -9.58
$ws.Range("Y4").Value = -9.039999999999999
$ws.Range("Z4").Value = -4.5
$ws.Range("AA4").Value = 92.06
$ws.Range("AB4").Value = 740.2
$ws.Range("AC4").Value = -352
$ws.Range("AD4").Value = -16.88
$ws.Range("AE4").Value = 4333
$ws.Range("AF4").Value = 1.37
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 52984990

# Row 5 - 2017/12 (IFRS연결)
$ws.Range("D5").Value = 1850
$ws.Range("E5").Value = 90
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 21
$ws.Range("I5").Value = 21
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3906
$ws.Range("L5").Value = 1755
$ws.Range("M5").Value = 2151
$ws.Range("N5").Value = 2148
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 276
$ws.Range("Q5").Value = 101
$ws.Range("R5").Value = 20
$ws.Range("S5").Value = -85
$ws.Range("T5").Value = 31
$ws.Range("U5").Value = 71
$ws.Range("V5").Value = 1126
$ws.Range("W5").Value = 4.88
$ws.Range("X5").Value = 1.15
$ws.Range("Y5").Value = 0.98
$ws.Range("Z5").Value = 0.53
$ws.Range("AA5").Value = 81.58
$ws.Range("AB5").Value = 746.05
$ws.Range("AC5").Value = 38
$ws.Range("AD5").Value = 200.54
$ws.Range("AE5").Value = 4299
$ws.Range("AF5").Value = 1.79
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 52984990
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6 - 2018/12 (IFRS연결)
$ws.Range("D6").Value = 1874
$ws.Range("E6").Value = 69
$ws.Range("F6").Value = 69
$ws.Range("G6").Value = 43
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 20
$ws.Range("K6").Value = 3851
$ws.Range("L6").Value = 1735
$ws.Range("M6").Value = 2117
$ws.Range("N6").Value = 2106
$ws.Range("P6").Value = 276
$ws.Range("Q6").Value = 321
$ws.Range("R6").Value = -103
$ws.Range("S6").Value = -42
$ws.Range("T6").Value = 18
$ws.Range("U6").Value = 303
$ws.Range("V6").Value = 1074
$ws.Range("W6").Value = 3.69
$ws.Range("X6").Value = 1.05
$ws.Range("Y6").Value = 0.92
$ws.Range("Z6").Value = 0.51
$ws.Range("AA6").Value = 81.95999999999999
$ws.Range("AB6").Value = 744.6799999999999
$ws.Range("AC6").Value = 36
$ws.Range("AD6").Value = 180.16
$ws.Range("AE6").Value = 4215
$ws.Range("AF6").Value = 1.52
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 52984990
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Rows 7-9 - 2019/12(E), 2020/12(E), 2021/12(E): clear all figures (D:AI),
# keep only the index/label columns A:C
$ws.Range("D7:AI9").ClearContents()

